$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicated/incorrect "Number" mapping values in columns E and G for
# the grid references D0-D6, E0-E6, F0-F6, G0-G6 (rows 24-51) so that each
# row of the ship-location lookup table has a unique number.
$newValues = @{
    24 = 30; 25 = 31; 26 = 32; 27 = 33; 28 = 34; 29 = 35; 30 = 36;
    31 = 40; 32 = 41; 33 = 42; 34 = 43; 35 = 44; 36 = 45; 37 = 46;
    38 = 50; 39 = 51; 40 = 52; 41 = 53; 42 = 54; 43 = 55; 44 = 56;
    45 = 60; 46 = 61; 47 = 62; 48 = 63; 49 = 64; 50 = 65; 51 = 66
}

foreach ($row in $newValues.Keys) {
    $val = $newValues[$row]
    $ws.Cells.Item($row, 5).Value = $val   # column E
    $ws.Cells.Item($row, 7).Value = $val   # column G
}

# Update the "Enter Ship Location Value Below" input cell (B3) with a new
# sample value, which drives the VLOOKUP result shown in B6.
$ws.Range("B3").Value = 63

# Restore sensible (non-zero) widths for the helper columns H:J; they remain
# hidden, but are no longer zero-width.
$ws.Columns("H:J").ColumnWidth = 10
$ws.Columns("H:J").EntireColumn.Hidden = $true

# Update the selected cell shown when the workbook is opened.
$ws.Range("B6").Select()
